$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1200
$ws.Range("J43").Value = 900
$ws.Range("L43").Value = 900
$ws.Range("N43").Value = -1038
$ws.Range("H58").Value = 2555.4443
$ws.Range("J58").Value = 2202.8
$ws.Range("L58").Value = 6608.400000000001
$ws.Range("N58").Value = -6908.400000000001
$ws.Range("H62").Value = 6676.1055
$ws.Range("I62").Value = 5824.467
$ws.Range("K62").Value = 5824.467
$ws.Range("M62").Value = -5200.467
$ws.Range("H65").Value = 6676.1055
$ws.Range("I65").Value = 5824.467
$ws.Range("K65").Value = 29122.335
$ws.Range("M65").Value = -26002.335
$ws.Range("H112").Value = 796956.1
$ws.Range("I112").Value = 1632.6666
$ws.Range("J112").Value = 858134.9
$ws.Range("K112").Value = 4897.9998
$ws.Range("L112").Value = 2574404.7
$ws.Range("M112").Value = -3789.9998
$ws.Range("N112").Value = -2576620.7
$ws.Range("H119").Value = 3497.5
$ws.Range("J119").Value = 2000
$ws.Range("L119").Value = 6000
$ws.Range("N119").Value = -15676
$ws.Range("H132").Value = 5654.8613
$ws.Range("I132").Value = 6393.037
$ws.Range("K132").Value = 19179.111
$ws.Range("M132").Value = -16649.111
$ws.Range("H135").Value = 947.913
$ws.Range("I135").Value = 809.1818
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 7282.6362
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -4747.6362
$ws.Range("N135").Value = -41070
$ws.Range("H137").Value = 31006.25
$ws.Range("I137").Value = 47128.89
$ws.Range("K137").Value = 141386.67
$ws.Range("M137").Value = -138836.67
$ws.Range("H138").Value = 3366.7144
$ws.Range("I138").Value = 1914.1389
$ws.Range("J138").Value = 5981.35
$ws.Range("K138").Value = 5742.4167
$ws.Range("L138").Value = 17944.05
$ws.Range("M138").Value = -602.4166999999998
$ws.Range("N138").Value = -28224.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13237690
$ws.Range("I32").Value = 6667982
$ws.Range("K32").Value = 6667982
$ws.Range("M32").Value = -6667695
$ws.Range("H61").Value = 3627.7144
$ws.Range("I61").Value = 3327
$ws.Range("J61").Value = 4590
$ws.Range("K61").Value = 3327
$ws.Range("L61").Value = 4590
$ws.Range("M61").Value = -3115
$ws.Range("N61").Value = -5014
$ws.Range("H63").Value = 159033.89
$ws.Range("I63").Value = 2376.25
$ws.Range("J63").Value = 284360
$ws.Range("K63").Value = 2376.25
$ws.Range("L63").Value = 284360
$ws.Range("M63").Value = -1690.25
$ws.Range("N63").Value = -285732
$ws.Range("H66").Value = 159033.89
$ws.Range("I66").Value = 2376.25
$ws.Range("J66").Value = 284360
$ws.Range("K66").Value = 11881.25
$ws.Range("L66").Value = 1421800
$ws.Range("M66").Value = -8449.25
$ws.Range("N66").Value = -1428664
$ws.Range("H74").Value = 3299.375
$ws.Range("I74").Value = 3522.75
$ws.Range("J74").Value = 3076
$ws.Range("K74").Value = 3522.75
$ws.Range("L74").Value = 3076
$ws.Range("M74").Value = -2648.75
$ws.Range("N74").Value = -4824
$ws.Range("H77").Value = 3299.375
$ws.Range("I77").Value = 3522.75
$ws.Range("J77").Value = 3076
$ws.Range("K77").Value = 17613.75
$ws.Range("L77").Value = 15380
$ws.Range("M77").Value = -13245.75
$ws.Range("N77").Value = -24116
$ws.Range("H124").Value = 25499.5
$ws.Range("J124").Value = 25499.5
$ws.Range("L124").Value = 25499.5
$ws.Range("N124").Value = -35319.5
$ws.Range("H132").Value = 3320
$ws.Range("I132").Value = 2654.8125
$ws.Range("K132").Value = 7964.4375
$ws.Range("M132").Value = -5434.4375
$ws.Range("H136").Value = 3627.7144
$ws.Range("I136").Value = 3327
$ws.Range("J136").Value = 4590
$ws.Range("K136").Value = 9981
$ws.Range("L136").Value = 13770
$ws.Range("M136").Value = -7431
$ws.Range("N136").Value = -18870

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 52355.25
$ws.Range("J2").Value = 52355.25
$ws.Range("L2").Value = 52355.25
$ws.Range("N2").Value = -52581.25
$ws.Range("H134").Value = 20240436
$ws.Range("I134").Value = 5496437
$ws.Range("J134").Value = 47622148
$ws.Range("K134").Value = 16489311
$ws.Range("L134").Value = 142866444
$ws.Range("M134").Value = -16486776
$ws.Range("N134").Value = -142871514

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 180.53334
$ws.Range("I7").Value = 71.125
$ws.Range("K7").Value = 71.125
$ws.Range("M7").Value = 41.875
$ws.Range("H22").Value = 1352.2
$ws.Range("I22").Value = 1647.625
$ws.Range("J22").Value = 170.5
$ws.Range("K22").Value = 1647.625
$ws.Range("L22").Value = 170.5
$ws.Range("M22").Value = -1297.625
$ws.Range("N22").Value = -870.5
$ws.Range("H31").Value = 4894.0312
$ws.Range("I31").Value = 4460.85
$ws.Range("K31").Value = 4460.85
$ws.Range("M31").Value = -4165.85
$ws.Range("H34").Value = 4894.0312
$ws.Range("I34").Value = 4460.85
$ws.Range("K34").Value = 4460.85
$ws.Range("M34").Value = -4258.85
$ws.Range("H50").Value = 39997.8
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 39997.8
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 39997.8
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -41247.8
$ws.Range("H58").Value = 3674.05
$ws.Range("I58").Value = 2567.75
$ws.Range("K58").Value = 2567.75
$ws.Range("M58").Value = -2364.75
$ws.Range("H132").Value = 2125
$ws.Range("I132").Value = 2125
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6375
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3845
$ws.Range("N132").ClearContents()
$ws.Range("H134").Value = 3251.4285
$ws.Range("I134").Value = 2683.7273
$ws.Range("K134").Value = 8051.1819
$ws.Range("M134").Value = -5516.1819
$ws.Range("H135").Value = 57199.8
$ws.Range("J135").Value = 57199.8
$ws.Range("L135").Value = 57199.8
$ws.Range("N135").Value = -67339.8
$ws.Range("H136").Value = 3674.05
$ws.Range("I136").Value = 2567.75
$ws.Range("K136").Value = 7703.25
$ws.Range("M136").Value = -5153.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 5024
$ws.Range("I94").Value = 5024
$ws.Range("K94").Value = 15072
$ws.Range("M94").Value = -14396
$ws.Range("H122").Value = 1968.7778
$ws.Range("J122").Value = 2089.875
$ws.Range("L122").Value = 18808.875
$ws.Range("N122").Value = -23708.875
$ws.Range("H129").Value = 1056.1765
$ws.Range("I129").Value = 399.35715
$ws.Range("J129").Value = 4121.3335
$ws.Range("K129").Value = 1198.07145
$ws.Range("L129").Value = 12364.0005
$ws.Range("M129").Value = 3801.92855
$ws.Range("N129").Value = -22364.0005
$ws.Range("H137").Value = 4428.6
$ws.Range("J137").Value = 2660.75
$ws.Range("L137").Value = 7982.25
$ws.Range("N137").Value = -18182.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4458.892
$ws.Range("I70").Value = 3990
$ws.Range("K70").Value = 3990
$ws.Range("M70").Value = -3720
$ws.Range("H73").Value = 4458.892
$ws.Range("I73").Value = 3990
$ws.Range("K73").Value = 3990
$ws.Range("M73").Value = -3054
$ws.Range("H80").Value = 2979.625
$ws.Range("J80").Value = 2979.625
$ws.Range("L80").Value = 2979.625
$ws.Range("N80").Value = -4975.625
$ws.Range("H83").Value = 2979.625
$ws.Range("J83").Value = 2979.625
$ws.Range("L83").Value = 14898.125
$ws.Range("N83").Value = -24882.125
$ws.Range("H122").Value = 7497.8335
$ws.Range("I122").Value = 7597.6
$ws.Range("J122").Value = 6999
$ws.Range("K122").Value = 22792.8
$ws.Range("L122").Value = 20997
$ws.Range("M122").Value = -20342.8
$ws.Range("N122").Value = -25897
$ws.Range("H132").Value = 3668.25
$ws.Range("I132").Value = 3446.5557
$ws.Range("K132").Value = 10339.6671
$ws.Range("M132").Value = -7809.667099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4182.875
$ws.Range("I40").Value = 4199.727
$ws.Range("J40").Value = 3997.5
$ws.Range("K40").Value = 4199.727
$ws.Range("L40").Value = 3997.5
$ws.Range("M40").Value = -4063.727
$ws.Range("N40").Value = -4269.5
$ws.Range("H136").Value = 4911.3335
$ws.Range("I136").Value = 4282
$ws.Range("K136").Value = 12846
$ws.Range("L136").Value = 17802
$ws.Range("M136").Value = -10296
$ws.Range("N136").Value = -22902

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4968.0835
$ws.Range("J62").Value = 5848.5
$ws.Range("L62").Value = 5848.5
$ws.Range("N62").Value = -7096.5
$ws.Range("H65").Value = 4968.0835
$ws.Range("J65").Value = 5848.5
$ws.Range("L65").Value = 29242.5
$ws.Range("N65").Value = -35482.5
$ws.Range("H132").Value = 3130.361
$ws.Range("I132").Value = 2636.2593
$ws.Range("K132").Value = 7908.777900000001
$ws.Range("M132").Value = -5378.777900000001
$ws.Range("H135").Value = 97378.336
$ws.Range("J135").Value = 97378.336
$ws.Range("L135").Value = 97378.336
$ws.Range("N135").Value = -107518.336
